$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '26.429.73'
$ws.Range("E2").Value = '  +1.29%  '

# Row 3
$ws.Range("D3").Value = '1.692.52'
$ws.Range("E3").Value = '  +1.45%  '

# Row 4
$ws.Range("E4").Value = '  +0.87%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '219.01'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.35%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5517'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +8.07%  '

# Row 8
$ws.Range("E8").Value = '  +1.58%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06480'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.55%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.11'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.08%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07703'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.47%  '

# Row 12
$ws.Range("D12").Value = '1.689.65'
$ws.Range("E12").Value = '  +0.65%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.549'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.76%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5825'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.19%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.000008423'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.91%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.28'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.79%  '

# Row 17
$ws.Range("D17").Value = '26.522.39'
$ws.Range("E17").Value = '  +2.47%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.954'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.43%  '

# Row 19
$ws.Range("E19").Value = '  +0.83%  '

# Row 20
$ws.Range("E20").Value = '  +1.67%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '190.28'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.18%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.240'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.75%  '

# Row 23
$ws.Range("E23").Value = '  +0.85%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '150.06'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.74%  '

# Row 25
$ws.Range("E25").Value = '  +6.96%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.901'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.98%  '

# Row 27
$ws.Range("E27").Value = '  +0.14%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.427'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +6.19%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.06324'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -5.24%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.330'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.17%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.588'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.57%  '

# Row 32
$ws.Range("E32").Value = '  +2.03%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.677'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.69%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.043'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.43%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6211'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.68%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.409'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.87%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.727'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.53%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.231'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.61%  '

# Row 39
$ws.Range("D39").Value = '1.124.24'
$ws.Range("E39").Value = '  +2.56%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01641'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.85%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8823'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.91%  '

# Row 42
$ws.Range("E42").Value = '  +0.80%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '101.01'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.38%  '

# Row 44
$ws.Range("D44").Value = '1.842.93'
$ws.Range("E44").Value = '  +1.53%  '

# Row 45
$ws.Range("B45").Value = 'BabyDogeCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00000000110'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.25%  '

# Row 46
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '57.39'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.81%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.232'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.47%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.008'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.43%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05285'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.01%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4305'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.66%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.077'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.33%  '
